$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Conflicting Features"
$ws.Range("G1").Value = "Third Party Domains"

$ws.Range("A2").Value = "https://www.google.com/"
$ws.Range("C2").Value = $false
$ws.Range("D2").Value = $false
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = "[]"
$ws.Range("G2").Value = "[]"
